# Lot of updates since last push
# Prefix each strain name in the menu with its sequential number ("N: Name"),
# numbering down column B (rows 4-23) first, then column H (rows 4-20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$num = 1

for ($row = 4; $row -le 23; $row++) {
    $cell = $ws.Range("B$row")
    $name = $cell.Value2
    if ($name -ne $null -and $name -ne "") {
        $cell.Value = "$($num): $name"
        $num++
    }
}

for ($row = 4; $row -le 20; $row++) {
    $cell = $ws.Range("H$row")
    $name = $cell.Value2
    if ($name -ne $null -and $name -ne "") {
        $cell.Value = "$($num): $name"
        $num++
    }
}
